# Fruta / hortaliza, semanal
# Insert two new weekly price records for "Ají" (Macroferia Regional de Talca)
# right after the existing row 334, shifting all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 335-336 (pushes old 335.. down to 337..)
$ws.Range("A335:A336").EntireRow.Insert()

# --- New row 335: Ají / Americana (o) ---
$ws.Cells.Item(335, 1).Value = 5
$ws.Cells.Item(335, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(335, 3).Value = "Maule"
$ws.Cells.Item(335, 4).Value = 45021
$ws.Cells.Item(335, 5).Value = 7
$ws.Cells.Item(335, 6).Value = 100112021
$ws.Cells.Item(335, 7).Value = "Ají"
$ws.Cells.Item(335, 8).Value = "Americana (o)"
$ws.Cells.Item(335, 9).Value = "Primera"
$ws.Cells.Item(335, 10).Value = 150
$ws.Cells.Item(335, 11).Value = 8000
$ws.Cells.Item(335, 12).Value = 8000
$ws.Cells.Item(335, 13).Value = 8000
$ws.Cells.Item(335, 14).Value = "$/caja 14 kilos"
$ws.Cells.Item(335, 15).Value = "Región del Maule"
$ws.Cells.Item(335, 16).Value = 571
$ws.Cells.Item(335, 17).Value = 14
$ws.Cells.Item(335, 18).Value = "Hortaliza"

# --- New row 336: Ají / Cacho cabra verde ---
$ws.Cells.Item(336, 1).Value = 5
$ws.Cells.Item(336, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(336, 3).Value = "Maule"
$ws.Cells.Item(336, 4).Value = 45021
$ws.Cells.Item(336, 5).Value = 7
$ws.Cells.Item(336, 6).Value = 100112021
$ws.Cells.Item(336, 7).Value = "Ají"
$ws.Cells.Item(336, 8).Value = "Cacho cabra verde"
$ws.Cells.Item(336, 9).Value = "Primera"
$ws.Cells.Item(336, 10).Value = 150
$ws.Cells.Item(336, 11).Value = 8000
$ws.Cells.Item(336, 12).Value = 8000
$ws.Cells.Item(336, 13).Value = 8000
$ws.Cells.Item(336, 14).Value = "$/caja 14 kilos"
$ws.Cells.Item(336, 15).Value = "Región del Maule"
$ws.Cells.Item(336, 16).Value = 571
$ws.Cells.Item(336, 17).Value = 14
$ws.Cells.Item(336, 18).Value = "Hortaliza"
